# Update automàtic: dades i banners [2026-02-10 18:20]
# Refreshes DATA_EXTRACCIO timestamps and several measurement columns
# (HUMITAT_MITJANA_DIA, PRECIPITACIO_ACUM_DIA, PRESSIO_ATMOSFERICA,
# RADIACIO_GLOBAL, RATXA_VENT_MAX, TEMPERATURA_MAXIMA_DIA,
# TEMPERATURA_MITJANA_DIA) with the latest meteocat readings.
# Values are prefixed with a leading apostrophe so Excel keeps them as
# plain text (matching the original inlineStr/General formatting)
# instead of auto-converting percentages or dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'2026-02-10 18:18:39"
$ws.Range("I2").Value = "'34.0 mm"
$ws.Range("O2").Value = "'0.8 °C"
$ws.Range("E3").Value = "'2026-02-10 18:18:42"
$ws.Range("I3").Value = "'21.8 mm"
$ws.Range("E4").Value = "'2026-02-10 18:18:44"
$ws.Range("J4").Value = "'1004.1 hPa"
$ws.Range("L4").Value = "'21.2 km/h - 90º 17:42 TU"
$ws.Range("O4").Value = "'11.9 °C"
$ws.Range("E5").Value = "'2026-02-10 18:18:47"
$ws.Range("I5").Value = "'28.9 mm"
$ws.Range("E6").Value = "'2026-02-10 18:18:49"
$ws.Range("O6").Value = "'9.6 °C"
$ws.Range("E7").Value = "'2026-02-10 18:18:51"
$ws.Range("O7").Value = "'15.0 °C"
$ws.Range("E8").Value = "'2026-02-10 18:18:54"
$ws.Range("H8").Value = "'85%"
$ws.Range("O8").Value = "'11.7 °C"
$ws.Range("E9").Value = "'2026-02-10 18:18:57"
$ws.Range("O9").Value = "'8.8 °C"
$ws.Range("E10").Value = "'2026-02-10 18:19:00"
$ws.Range("E11").Value = "'2026-02-10 18:19:02"
$ws.Range("O11").Value = "'7.5 °C"
$ws.Range("E12").Value = "'2026-02-10 18:19:05"
$ws.Range("O12").Value = "'8.9 °C"
$ws.Range("E13").Value = "'2026-02-10 18:19:07"
$ws.Range("I13").Value = "'3.4 mm"
$ws.Range("J13").Value = "'1006.9 hPa"
$ws.Range("O13").Value = "'5.0 °C"
$ws.Range("E14").Value = "'2026-02-10 18:19:10"
$ws.Range("E15").Value = "'2026-02-10 18:19:13"
$ws.Range("E16").Value = "'2026-02-10 18:19:15"
$ws.Range("I16").Value = "'22.5 mm"
$ws.Range("E17").Value = "'2026-02-10 18:19:17"
$ws.Range("L17").Value = "'120.6 km/h - 270º 17:55 TU"
$ws.Range("E18").Value = "'2026-02-10 18:19:20"
$ws.Range("J18").Value = "'1004.5 hPa"
$ws.Range("E19").Value = "'2026-02-10 18:19:22"
$ws.Range("L19").Value = "'32.4 km/h - 228º 17:31 TU"
$ws.Range("O19").Value = "'6.3 °C"
$ws.Range("E20").Value = "'2026-02-10 18:19:25"
$ws.Range("I20").Value = "'4.7 mm"
$ws.Range("L20").Value = "'62.3 km/h - 296º 17:58 TU"
$ws.Range("M20").Value = "'2.6 °C 17:53 TU"
$ws.Range("E21").Value = "'2026-02-10 18:19:28"
$ws.Range("I21").Value = "'4.5 mm"
$ws.Range("J21").Value = "'1006.3 hPa"
$ws.Range("O21").Value = "'7.1 °C"
$ws.Range("E22").Value = "'2026-02-10 18:19:30"
$ws.Range("I22").Value = "'8.0 mm"
$ws.Range("O22").Value = "'-0.6 °C"
$ws.Range("E23").Value = "'2026-02-10 18:19:33"
$ws.Range("I23").Value = "'22.7 mm"
$ws.Range("E24").Value = "'2026-02-10 18:19:36"
$ws.Range("J24").Value = "'1006.2 hPa"
$ws.Range("O24").Value = "'11.2 °C"
$ws.Range("E25").Value = "'2026-02-10 18:19:39"
$ws.Range("I25").Value = "'13.5 mm"
$ws.Range("E26").Value = "'2026-02-10 18:19:41"
$ws.Range("E27").Value = "'2026-02-10 18:19:44"
$ws.Range("I27").Value = "'2.5 mm"
$ws.Range("L27").Value = "'48.2 km/h - 228º 17:34 TU"
$ws.Range("E28").Value = "'2026-02-10 18:19:47"
$ws.Range("J28").Value = "'1004.4 hPa"
$ws.Range("O28").Value = "'8.8 °C"
$ws.Range("E29").Value = "'2026-02-10 18:19:49"
$ws.Range("O29").Value = "'10.5 °C"
$ws.Range("E30").Value = "'2026-02-10 18:19:52"
$ws.Range("E31").Value = "'2026-02-10 18:19:54"
$ws.Range("H31").Value = "'81%"
$ws.Range("E32").Value = "'2026-02-10 18:19:57"
$ws.Range("O32").Value = "'10.3 °C"
$ws.Range("E33").Value = "'2026-02-10 18:20:00"
$ws.Range("I33").Value = "'8.8 mm"
$ws.Range("J33").Value = "'1006.7 hPa"
$ws.Range("O33").Value = "'4.0 °C"
$ws.Range("E34").Value = "'2026-02-10 18:20:02"
$ws.Range("E35").Value = "'2026-02-10 18:20:05"
$ws.Range("J35").Value = "'1005.1 hPa"
$ws.Range("O35").Value = "'12.9 °C"
$ws.Range("E36").Value = "'2026-02-10 18:20:08"
$ws.Range("J36").Value = "'1004.6 hPa"
$ws.Range("E37").Value = "'2026-02-10 18:20:11"
$ws.Range("J37").Value = "'1005.8 hPa"
$ws.Range("O37").Value = "'6.5 °C"
$ws.Range("E38").Value = "'2026-02-10 18:20:13"
$ws.Range("H38").Value = "'88%"
$ws.Range("O38").Value = "'10.6 °C"
$ws.Range("E39").Value = "'2026-02-10 18:20:16"
$ws.Range("L39").Value = "'65.9 km/h - 287º 17:53 TU"
$ws.Range("E40").Value = "'2026-02-10 18:20:19"
$ws.Range("I40").Value = "'7.9 mm"
$ws.Range("J40").Value = "'1007.2 hPa"
$ws.Range("O40").Value = "'7.4 °C"
$ws.Range("E41").Value = "'2026-02-10 18:20:21"
$ws.Range("E42").Value = "'2026-02-10 18:20:24"
$ws.Range("E43").Value = "'2026-02-10 18:20:26"
$ws.Range("H43").Value = "'88%"
$ws.Range("O43").Value = "'9.2 °C"
$ws.Range("E44").Value = "'2026-02-10 18:20:29"
$ws.Range("I44").Value = "'21.1 mm"
$ws.Range("O44").Value = "'0.3 °C"
$ws.Range("E45").Value = "'2026-02-10 18:20:32"
$ws.Range("H45").Value = "'94%"
$ws.Range("I45").Value = "'27.8 mm"
$ws.Range("J45").Value = "'1005.8 hPa"
$ws.Range("O45").Value = "'6.4 °C"
$ws.Range("E46").Value = "'2026-02-10 18:20:35"
$ws.Range("J46").Value = "'1006.1 hPa"
$ws.Range("K46").Value = "'8.8 MJ/m2"
$ws.Range("O46").Value = "'14.0 °C"
Write-Host "Applied 111 cell updates to sheet1"
